# Update the cached "today" date shown by the datetimeFigureOut field
# placeholders on the slide master and on every slide layout
# (04.11.18 -> 11.11.18).
$p = $ppt.ActivePresentation

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "11.11.18"
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "11.11.18"
        }
    }
}

# Relabel the two latent-variable textboxes from z1/z2 to y1/y2 and widen
# them slightly so the new caption still fits ("spAutoFit" boxes).
$s = $p.Slides.Item(1)
$newWidth = 548548 / 12700

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "z1") {
            $shp.TextFrame.TextRange.Text = "y1"
            $shp.Width = $newWidth
        } elseif ($t -eq "z2") {
            $shp.TextFrame.TextRange.Text = "y2"
            $shp.Width = $newWidth
        }
    }
}
